$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.317023666666666
$ws.Range("H2").Value = 9.951070999999999
$ws.Range("I2").Value = 0.005526424869593284
$ws.Range("J2").Value = 0.005526424869593283
$ws.Range("M2").Value = 7.393757999999999
$ws.Range("N2").Value = 22.181274
$ws.Range("O2").Value = 0.6880101757622353
$ws.Range("P2").Value = 0.6880101757622353
$ws.Range("Q2").Value = 24.52527027160599
$ws.Range("R2").Value = 220.727432444454
$ws.Range("S2").Value = 0.003802236545865664
$ws.Range("T2").Value = 0.003802236545865663
$ws.Range("G3").Value = 3.317023666666666
$ws.Range("H3").Value = 9.951070999999999
$ws.Range("I3").Value = 0.005526424869593284
$ws.Range("J3").Value = 0.005526424869593283
$ws.Range("O3").Value = 0.2613149228497358
$ws.Range("P3").Value = 0.2613149228497358
$ws.Range("Q3").Value = 9.315006281401889
$ws.Range("R3").Value = 83.835056532617
$ws.Range("S3").Value = 0.00144413728843263
$ws.Range("T3").Value = 0.00144413728843263
$ws.Range("G4").Value = 3.317023666666666
$ws.Range("H4").Value = 9.951070999999999
$ws.Range("I4").Value = 0.005526424869593284
$ws.Range("J4").Value = 0.005526424869593283
$ws.Range("M4").Value = 0.544582
$ws.Range("N4").Value = 1.633746
$ws.Range("O4").Value = 0.05067490138802889
$ws.Range("P4").Value = 0.05067490138802888
$ws.Range("Q4").Value = 1.806391382440666
$ws.Range("R4").Value = 16.257522441966
$ws.Range("S4").Value = 0.0002800510352949901
$ws.Range("T4").Value = 0.00028005103529499
$ws.Range("I5").Value = 0.9618737623074323
$ws.Range("J5").Value = 0.961873762307432
$ws.Range("M5").Value = 7.393757999999999
$ws.Range("N5").Value = 22.181274
$ws.Range("O5").Value = 0.6880101757622353
$ws.Range("P5").Value = 0.6880101757622353
$ws.Range("Q5").Value = 4268.621132904752
$ws.Range("R5").Value = 38417.59019614277
$ws.Range("S5").Value = 0.6617789362662191
$ws.Range("T5").Value = 0.6617789362662189
$ws.Range("I6").Value = 0.9618737623074323
$ws.Range("J6").Value = 0.961873762307432
$ws.Range("O6").Value = 0.2613149228497358
$ws.Range("P6").Value = 0.2613149228497358
$ws.Range("S6").Value = 0.2513519679885518
$ws.Range("T6").Value = 0.2513519679885517
$ws.Range("I7").Value = 0.9618737623074323
$ws.Range("J7").Value = 0.961873762307432
$ws.Range("M7").Value = 0.544582
$ws.Range("N7").Value = 1.633746
$ws.Range("O7").Value = 0.05067490138802889
$ws.Range("P7").Value = 0.05067490138802888
$ws.Range("Q7").Value = 314.4022611775413
$ws.Range("R7").Value = 2829.620350597872
$ws.Range("S7").Value = 0.04874285805266147
$ws.Range("T7").Value = 0.04874285805266145
$ws.Range("G8").Value = 18.65467933333333
$ws.Range("H8").Value = 55.964038
$ws.Range("I8").Value = 0.03108017734031479
$ws.Range("J8").Value = 0.03108017734031478
$ws.Range("M8").Value = 7.393757999999999
$ws.Range("N8").Value = 22.181274
$ws.Range("O8").Value = 0.6880101757622353
$ws.Range("P8").Value = 0.6880101757622353
$ws.Range("Q8").Value = 137.928184558268
$ws.Range("R8").Value = 1241.353661024412
$ws.Range("S8").Value = 0.02138347827463142
$ws.Range("T8").Value = 0.02138347827463142
$ws.Range("G9").Value = 18.65467933333333
$ws.Range("H9").Value = 55.964038
$ws.Range("I9").Value = 0.03108017734031479
$ws.Range("J9").Value = 0.03108017734031478
$ws.Range("O9").Value = 0.2613149228497358
$ws.Range("P9").Value = 0.2613149228497358
$ws.Range("Q9").Value = 52.38686021862512
$ws.Range("R9").Value = 471.4817419676261
$ws.Range("S9").Value = 0.008121714143840468
$ws.Range("T9").Value = 0.008121714143840464
$ws.Range("G10").Value = 18.65467933333333
$ws.Range("H10").Value = 55.964038
$ws.Range("I10").Value = 0.03108017734031479
$ws.Range("J10").Value = 0.03108017734031478
$ws.Range("M10").Value = 0.544582
$ws.Range("N10").Value = 1.633746
$ws.Range("O10").Value = 0.05067490138802889
$ws.Range("P10").Value = 0.05067490138802888
$ws.Range("Q10").Value = 10.15900258070533
$ws.Range("R10").Value = 91.431023226348
$ws.Range("S10").Value = 0.001574984921842902
$ws.Range("T10").Value = 0.001574984921842901
$ws.Range("G11").Value = 0.9121026666666667
$ws.Range("H11").Value = 2.736308
$ws.Range("I11").Value = 0.001519635482659812
$ws.Range("J11").Value = 0.001519635482659812
$ws.Range("M11").Value = 7.393757999999999
$ws.Range("N11").Value = 22.181274
$ws.Range("O11").Value = 0.6880101757622353
$ws.Range("P11").Value = 0.6880101757622353
$ws.Range("Q11").Value = 6.743866388488
$ws.Range("R11").Value = 60.694797496392
$ws.Range("S11").Value = 0.001045524675519307
$ws.Range("T11").Value = 0.001045524675519307
$ws.Range("G12").Value = 0.9121026666666667
$ws.Range("H12").Value = 2.736308
$ws.Range("I12").Value = 0.001519635482659812
$ws.Range("J12").Value = 0.001519635482659812
$ws.Range("O12").Value = 0.2613149228497358
$ws.Range("P12").Value = 0.2613149228497358
$ws.Range("Q12").Value = 2.561405320879556
$ws.Range("R12").Value = 23.052647887916
$ws.Range("S12").Value = 0.0003971034289109699
$ws.Range("T12").Value = 0.0003971034289109698
$ws.Range("G13").Value = 0.9121026666666667
$ws.Range("H13").Value = 2.736308
$ws.Range("I13").Value = 0.001519635482659812
$ws.Range("J13").Value = 0.001519635482659812
$ws.Range("M13").Value = 0.544582
$ws.Range("N13").Value = 1.633746
$ws.Range("O13").Value = 0.05067490138802889
$ws.Range("P13").Value = 0.05067490138802888
$ws.Range("Q13").Value = 0.4967146944186667
$ws.Range("R13").Value = 4.470432249768
$ws.Range("S13").Value = 0.00007700737822953568
$ws.Range("T13").Value = 0.00007700737822953565
